$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "Country" column (column D), shifting Year/Quarter/Month/Day/Value left.
$ws.Range("D1").EntireColumn.Delete()

# Add the "2 X 3" testing-class rows: 2 extra years (2019/2018) already exist for
# Quarter 2 & 3 at Day=22; add a second Day (23) slice covering the same
# 2 (years) X 2 (quarters) block for each of the 4 companies -> 16 new rows.
$data = @(
  @("Healthcare","Medical Laboratories & Research","Agilent Technologies, Inc.",2019,2,6,23,74.34),
  @("Basic Materials","Aluminum","Alcoa Corporation",2019,2,6,23,23.6),
  @("Financial","Exchange Traded Fund","Perth Mint Physical Gold ETF",2019,2,6,23,14.98),
  @("Financial","Asset Management","Altaba Inc.",2019,2,6,23,69.88),
  @("Healthcare","Medical Laboratories & Research","Agilent Technologies, Inc.",2019,3,8,23,78.34),
  @("Basic Materials","Aluminum","Alcoa Corporation",2019,3,8,23,21.6),
  @("Financial","Exchange Traded Fund","Perth Mint Physical Gold ETF",2019,3,8,23,13.98),
  @("Financial","Asset Management","Altaba Inc.",2019,3,8,23,68.88),
  @("Healthcare","Medical Laboratories & Research","Agilent Technologies, Inc.",2018,2,6,23,73.34),
  @("Basic Materials","Aluminum","Alcoa Corporation",2018,2,6,23,22.6),
  @("Financial","Exchange Traded Fund","Perth Mint Physical Gold ETF",2018,2,6,23,13.98),
  @("Financial","Asset Management","Altaba Inc.",2018,2,6,23,64.88),
  @("Healthcare","Medical Laboratories & Research","Agilent Technologies, Inc.",2018,3,8,23,65.34),
  @("Basic Materials","Aluminum","Alcoa Corporation",2018,3,8,23,22.6),
  @("Financial","Exchange Traded Fund","Perth Mint Physical Gold ETF",2018,3,8,23,13.98),
  @("Financial","Asset Management","Altaba Inc.",2018,3,8,23,63.88)
)

$r = 18
foreach ($row in $data) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}

[void]$ws.Range("K12").Select()
